$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "last refreshed" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 21:16"

# Full refreshed & re-sorted (by "Casos totales" desc) country table for rows 4-187.
# Each entry: Pais, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$data = @(
    @('China',80967,39,71150,6569,2136,3,3248),
    @('Italia',47021,5986,5129,37860,2655,627,4032),
    @('España',20412,2335,1588,17774,939,219,1050),
    @('Alemania',19848,4528,180,19601,2,23,67),
    @('Iran',19644,1237,6745,11466,0,149,1433),
    @('Estados Unidos',18121,4332,125,17763,64,26,233),
    @('Francia',12612,1617,1295,10867,1122,78,450),
    @('Corea del Sur',8652,87,2233,6325,59,3,94),
    @('Suiza',5407,1185,15,5336,0,13,56),
    @('Reino Unido',3983,714,65,3741,20,33,177),
    @('Paises Bajos',2994,534,2,2886,210,30,106),
    @('Austria',2491,312,9,2476,14,0,6),
    @('Belgica',2257,462,204,2016,164,16,37),
    @('Noruega',1926,136,1,1918,27,0,7),
    @('Suecia',1639,200,16,1607,21,5,16),
    @('Dinamarca',1255,104,1,1245,37,3,9),
    @('Malasia',1030,130,87,940,26,1,3),
    @('Portugal',1020,234,5,1009,26,2,6),
    @('Japon',963,20,215,715,50,0,33),
    @('Canada',943,70,11,920,1,0,12),
    @('Australia',876,120,46,823,2,0,7),
    @('Chequia',833,139,4,829,6,0,0),
    @('Brasil',750,110,2,737,18,4,11),
    @('Crucero',712,0,527,177,14,1,8),
    @('Israel',705,28,15,690,10,0,0),
    @('Irlanda',683,126,5,675,6,0,3),
    @('Pakistan',500,46,13,484,0,1,3),
    @('Grecia',495,31,19,466,20,4,10),
    @('Luxemburgo',484,149,6,473,1,1,5),
    @('Catar',470,10,10,460,6,0,0),
    @('Finlandia',450,50,10,440,2,0,0),
    @('Chile',434,92,6,428,7,0,0),
    @('Polonia',425,70,13,407,3,0,5),
    @('Islandia',409,79,5,404,1,0,0),
    @('Singapur',385,40,131,254,14,0,0),
    @('Indonesia',369,60,17,320,0,7,32),
    @('Ecuador',367,107,3,359,2,2,5),
    @('Turquia',359,0,0,355,0,0,4),
    @('Arabia Saudita',344,70,8,336,0,0,0),
    @('Eslovenia',341,22,0,340,6,0,1),
    @('Tailandia',322,50,42,279,1,0,1),
    @('Rumania',308,31,31,277,11,0,0),
    @('Barein',297,18,125,171,4,0,1),
    @('Egipto',285,29,42,235,0,1,8),
    @('Estonia',283,16,1,282,1,0,0),
    @('Peru',263,29,1,258,5,3,4),
    @('Hong Kong',256,48,98,154,4,0,4),
    @('Rusia',253,54,12,240,0,0,1),
    @('India',249,55,23,221,0,1,5),
    @('Filipinas',230,13,8,204,1,1,18),
    @('Irak',208,16,49,142,0,4,17),
    @('Sudafrica',202,52,0,202,0,0,0),
    @('Libano',177,20,4,169,3,0,4),
    @('Mexico',164,46,4,159,1,0,1),
    @('Kuwait',159,11,22,137,5,0,0),
    @('Colombia',145,37,1,144,0,0,0),
    @('San Marino',144,0,4,126,12,0,14),
    @('Emiratos Arabes Unidos',140,0,31,109,2,0,0),
    @('Eslovaquia',137,13,0,137,2,0,0),
    @('Panama',137,0,1,135,7,0,1),
    @('Armenia',136,14,1,135,2,0,0),
    @('Serbia',135,32,2,132,4,1,1),
    @('Taiwan',135,27,28,105,0,1,2),
    @('Bulgaria',129,22,1,125,0,0,3),
    @('Argentina',128,0,3,122,0,0,3),
    @('Croacia',128,18,5,122,0,0,1),
    @('Costa Rica',113,26,2,109,2,1,2),
    @('Letonia',111,25,1,110,0,0,0),
    @('Uruguay',94,15,0,94,0,0,0),
    @('Vietnam',91,6,17,74,0,0,0),
    @('Argelia',90,0,32,47,0,2,11),
    @('Jordania',85,16,1,84,0,0,0),
    @('Hungria',85,12,7,74,6,3,4),
    @('Bosnia y Herzegovina',84,20,2,82,1,0,0),
    @('Islas Feroe',80,8,3,77,0,0,0),
    @('Brunei',78,5,1,77,2,0,0),
    @('Marruecos',77,14,2,72,1,1,3),
    @('Republica de Chipre',75,8,0,75,1,0,0),
    @('Principado de Andorra',75,1,1,74,2,0,0),
    @('Sri Lanka',73,13,3,70,0,0,0),
    @('Republica Dominicana',72,38,0,70,0,0,2),
    @('Republica de Macedonia',70,20,1,69,1,0,0),
    @('Albania',70,6,0,68,2,0,2),
    @('Bielorrusia',69,18,15,54,0,0,0),
    @('Moldavia',66,17,1,64,3,0,1),
    @('Malta',64,11,2,62,1,0,0),
    @('Lituania',63,15,1,62,1,0,0),
    @('Tunez',54,15,1,52,7,0,1),
    @('Kazajistan',52,8,0,52,0,0,0),
    @('Camboya',51,14,1,50,0,0,0),
    @('Oman',48,0,13,35,0,0,0),
    @('Estado de Palestina',48,1,17,31,0,0,0),
    @('Guadalupe',45,12,0,45,0,0,0),
    @('Georgia',44,4,1,43,1,0,0),
    @('Azerbaiyan',44,0,7,36,0,0,1),
    @('Venezuela',42,0,0,42,0,0,0),
    @('Burkina Faso',40,7,4,35,0,0,1),
    @('Nueva Zelanda',39,11,0,39,0,0,0),
    @('Ucrania',39,13,1,35,0,0,3),
    @('Senegal',38,2,2,36,0,0,0),
    @('Uzbekistan',33,10,0,33,0,0,0),
    @('Martinica',32,9,0,31,7,0,1),
    @('Liechtenstein',28,0,0,28,0,0,0),
    @('Reunion',28,0,0,28,0,0,0),
    @('Honduras',24,12,0,24,0,0,0),
    @('Afganistan',24,2,1,23,0,0,0),
    @('Camerun',20,7,2,18,0,0,0),
    @('Banglades',20,2,3,16,1,0,1),
    @('Consejo Danes para los Refugiados',18,4,0,18,0,0,0),
    @('Ruanda',17,6,0,17,0,0,0),
    @('Macao',17,0,10,7,0,0,0),
    @('Ghana',16,5,0,16,0,0,0),
    @('Bolivia',16,1,0,16,0,0,0),
    @('Cuba',16,5,0,15,0,0,1),
    @('Jamaica',16,1,2,13,0,0,1),
    @('Guayana Francesa',15,0,0,15,0,0,0),
    @('Guyana',15,10,0,14,0,0,1),
    @('Puerto Rico',14,8,0,14,0,0,0),
    @('Montenegro',14,1,0,14,0,0,0),
    @('Guam',14,2,0,14,0,0,0),
    @('Paraguay',13,0,0,13,1,0,0),
    @('Maldivas',13,0,0,13,0,0,0),
    @('Mauricio',12,5,0,12,0,0,0),
    @('Nigeria',12,0,1,11,0,0,0),
    @('Guatemala',12,3,0,11,0,0,1),
    @('Polinesia Francesa',11,5,0,11,0,0,0),
    @('Monaco',11,1,0,11,0,0,0),
    @('Gibraltar',10,0,2,8,0,0,0),
    @('Togo',9,8,0,9,0,0,0),
    @('Trinidad yTobago',9,0,0,9,0,0,0),
    @('Etiopia',9,2,0,9,0,0,0),
    @('Costa de Marfil',9,0,1,8,0,0,0),
    @('Kenia',7,0,0,7,0,0,0),
    @('Seychelles',7,1,0,7,0,0,0),
    @('Mongolia',6,0,0,6,0,0,0),
    @('Mayotte',6,2,0,6,0,0,0),
    @('Guinea Ecuatorial',6,0,0,6,0,0,0),
    @('Kirguistan',6,3,0,6,0,0,0),
    @('Tanzania',6,0,0,6,0,0,0),
    @('Barbados',5,0,0,5,0,0,0),
    @('Aruba',5,0,1,4,0,0,0),
    @('San Martin (Parte Francesa)',4,1,0,4,0,0,0),
    @('Surinam',4,3,0,4,0,0,0),
    @('Madagascar',3,3,0,3,0,0,0),
    @('San Bartolome',3,0,0,3,0,0,0),
    @('Republica de Africa Central',3,2,0,3,0,0,0),
    @('Bahamas',3,0,0,3,0,0,0),
    @('Congo',3,0,0,3,0,0,0),
    @('Islas Virgenes de los Estados Unidos',3,0,0,3,0,0,0),
    @('Namibia',3,0,0,3,0,0,0),
    @('Gabon',3,0,0,2,0,1,1),
    @('Islas Caimanes',3,0,0,2,0,0,1),
    @('Curazao',3,0,0,2,0,0,1),
    @('Isla de Man',2,1,0,2,0,0,0),
    @('Liberia',2,0,0,2,0,0,0),
    @('Butan',2,1,0,2,0,0,0),
    @('Zambia',2,0,0,2,0,0,0),
    @('Groenlandia',2,0,0,2,0,0,0),
    @('Benin',2,0,0,2,0,0,0),
    @('Santa Lucia',2,0,0,2,0,0,0),
    @('Guinea',2,1,0,2,0,0,0),
    @('Mauritania',2,0,0,2,0,0,0),
    @('Nueva Caledonia',2,0,0,2,0,0,0),
    @('Bermudas',2,0,0,2,0,0,0),
    @('Haiti',2,2,0,2,0,0,0),
    @('Sudan',2,0,0,1,0,0,1),
    @('Santa Sede',1,0,0,1,0,0,0),
    @('Nicaragua',1,0,0,1,0,0,0),
    @('San Martin (Parte Holandesa)',1,0,0,1,0,0,0),
    @('Republica del Chad',1,0,0,1,0,0,0),
    @('Niger',1,0,0,1,0,0,0),
    @('Suazilandia',1,0,0,1,0,0,0),
    @('San Vicente y las Granadinas',1,0,0,1,0,0,0),
    @('Republica de Yibuti',1,0,0,1,0,0,0),
    @('Montserrat',1,0,0,1,0,0,0),
    @('Gambia',1,0,0,1,0,0,0),
    @('Fiyi',1,0,0,1,0,0,0),
    @('Cabo Verde',1,1,0,1,0,0,0),
    @('Somalia',1,0,0,1,0,0,0),
    @('El Salvador',1,0,0,1,0,0,0),
    @('Angola',1,1,0,1,0,0,0),
    @('Antigua y Barbuda',1,0,0,1,0,0,0),
    @('Papua Nueva Guinea',1,1,0,1,0,0,0),
    @('Nepal',1,0,1,0,0,0,0)
)

$r = 4
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r++
}
